# Move borehole-level funding to borehole.funding
# Adds a new "funding" column (Q) to the "borehole" sheet, with header,
# header style matching the other header cells, a descriptive cell
# comment (matching the style used by the other column comments), an
# appropriately widened column, and conditional-formatting formulas on
# the existing columns updated to include the new column in their
# "all columns blank" checks.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("borehole")

# --- 1. Header cell value ------------------------------------------------
$ws.Range("Q1").Value = "funding"

# --- 2. Header style (match the bold/shaded style used by A1:P1) --------
$ws.Range("Q1").Font.Bold = $true
$ws.Range("Q1").Interior.Color = $ws.Range("P1").Interior.Color

# --- 3. Column width (match the neighboring data columns) ---------------
$ws.Range("Q1").ColumnWidth = $ws.Range("P1").ColumnWidth

# --- 4. Cell comment describing the field --------------------------------
$commentText = "[string] funding" + [char]10 + "Funding sources as a pipe-delimited list. Each entry must be in the format {funder} [{rorid}] > {award} [{number}] ({url}), where only the funder is required, funder and award cannot contain parentheses, and rorid is the funder's ROR (https://ror.org) ID (e.g. 01jtrvx49)." + [char]10 + "constraints:" + [char]10 + "  - pattern: [^\s]+( [^\s]+)*"
$comment = $ws.Range("Q1").AddComment()
$comment.Text($commentText)

# --- 5. Extend "row fully blank" conditional-formatting formulas --------
# These formulas reference the full data-column span ($A2:$P2, 16 cols);
# now that there are 17 data columns ($A2:$Q2), update both the range and
# the blank-count they compare against.
$ranges = @("A2:A1048576", "B2:B1048576", "D2:D1048576", "E2:E1048576", "F2:F1048576")
foreach ($r in $ranges) {
    $fc = $ws.Range($r).FormatConditions.Item(1)
    $f = $fc.Formula1
    $f = $f.Replace("`$A2:`$P2", "`$A2:`$Q2")
    $f = $f.Replace("<> 16", "<> 17")
    $fc.Formula1 = $f
}
